# Generate Report for Handoff
# Update localization status from "In Translation" -> "Ready for handoff"
# and bump the related timestamps, matching a refreshed CI report run.

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn     = $wb.Worksheets.Item("zh-cn")
$wsDeDe     = $wb.Worksheets.Item("de-de")

# --- Overview sheet (row 2): zh-cn / de-de status columns + HO xliff date ---
$wsOverview.Range("E2").Value = "Ready for handoff"
$wsOverview.Range("F2").Value = "Ready for handoff"
$wsOverview.Range("G2").Value = "2016-09-02 21:07:51"

# --- zh-cn sheet (row 2): Status + Latest Handoff Datetime ---
$wsZhCn.Range("C2").Value = "Ready for handoff"
$wsZhCn.Range("H2").Value = "2016-09-02 21:07:46"

# --- de-de sheet (row 2): Status + Latest Handoff Datetime ---
$wsDeDe.Range("C2").Value = "Ready for handoff"
$wsDeDe.Range("H2").Value = "2016-09-02 21:07:51"

# The "Status" column text grew ("In Translation" -> "Ready for handoff"),
# so the report generator re-autofits that column's width on each sheet
# (target rendered width ~17.216 characters).
$newStatusColWidth = 16.333333333333332
$wsOverview.Columns.Item(5).ColumnWidth = $newStatusColWidth
$wsOverview.Columns.Item(6).ColumnWidth = $newStatusColWidth
$wsZhCn.Columns.Item(3).ColumnWidth = $newStatusColWidth
$wsDeDe.Columns.Item(3).ColumnWidth = $newStatusColWidth
